$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift "Calibrated By" from column I to new column J,
#     and give column I a new header "Config File Updated?"
$ws.Range("J1").Value = $ws.Range("I1").Value2
$ws.Range("I1").Value = "Config File Updated?"

# --- Row 2: shift existing "NA" from I2 to new J2, put "Y" in I2
$ws.Range("J2").Value = $ws.Range("I2").Value2
$ws.Range("I2").Value = "Y"

# --- Row 3: Calibration Power changes from "100mW " to "10mW "
$ws.Range("B3").Value = "10mW "
# shift existing "NA" from I3 to new J3, put "Y" in I3
$ws.Range("J3").Value = $ws.Range("I3").Value2
$ws.Range("I3").Value = "Y"

# --- New Row 4: new calibration entry
$ws.Range("A4").Value = 45890
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat

$ws.Range("B4").Value = "100mW "

$ws.Range("C4").Value = 910
$ws.Range("D4").Value = 8210
$ws.Range("E4").Formula = "=D4/C4"

$ws.Range("F4").Value = 30000
$ws.Range("F4").NumberFormat = "#,##0"

$ws.Range("G4").Value = 30.3
$ws.Range("H4").Formula = "=F4/G4"

$ws.Range("I4").Value = "Y"
$ws.Range("J4").Value = "NA"

# --- Column width adjustments
$ws.Columns.Item(6).ColumnWidth = 23.14    # column F -> width 24
$ws.Columns.Item(9).ColumnWidth = 37.14    # column I -> width 38
$ws.Columns.Item(10).ColumnWidth = 14.14   # column J (new) -> width 15

# --- Selection moved as part of editing session
$ws.Range("F14").Select()
